$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.68414866666667
$ws.Range("H2").Value = 41.052446
$ws.Range("I2").Value = 0.06687192512742757
$ws.Range("J2").Value = 0.06728034099283527
$ws.Range("M2").Value = 3.029459000000001
$ws.Range("N2").Value = 9.088377000000001
$ws.Range("O2").Value = 0.009650183621689718
$ws.Range("P2").Value = 0.0101382356457613
$ws.Range("Q2").Value = 41.45556733557134
$ws.Range("R2").Value = 373.100106020142
$ws.Range("S2").Value = 0.0006453263566155627
$ws.Range("T2").Value = 0.0006821039513125376
$ws.Range("G3").Value = 13.68414866666667
$ws.Range("H3").Value = 41.052446
$ws.Range("I3").Value = 0.06687192512742757
$ws.Range("J3").Value = 0.06728034099283527
$ws.Range("O3").Value = 0.2405828560439687
$ws.Range("P3").Value = 0.2527501840920355
$ws.Range("Q3").Value = 1033.503524855051
$ws.Range("R3").Value = 9301.531723695462
$ws.Range("S3").Value = 0.01608823873631496
$ws.Range("T3").Value = 0.01700511857171404
$ws.Range("G4").Value = 13.68414866666667
$ws.Range("H4").Value = 41.052446
$ws.Range("I4").Value = 0.06687192512742757
$ws.Range("J4").Value = 0.06728034099283527
$ws.Range("M4").Value = 109.55234
$ws.Range("N4").Value = 328.65702
$ws.Range("O4").Value = 0.348973264594696
$ws.Range("P4").Value = 0.3666223700220274
$ws.Range("Q4").Value = 1499.130507341213
$ws.Range("R4").Value = 13492.17456607092
$ws.Range("S4").Value = 0.02333651402145048
$ws.Range("T4").Value = 0.02466647807068343
$ws.Range("G5").Value = 13.68414866666667
$ws.Range("H5").Value = 41.052446
$ws.Range("I5").Value = 0.06687192512742757
$ws.Range("J5").Value = 0.06728034099283527
$ws.Range("M5").Value = 45.33717900000001
$ws.Range("N5").Value = 90.67435800000001
$ws.Range("O5").Value = 0.144419218824026
$ws.Range("P5").Value = 0.1011487538899543
$ws.Range("Q5").Value = 620.4006975632781
$ws.Range("R5").Value = 3722.404185379668
$ws.Range("S5").Value = 0.009657591188161847
$ws.Range("T5").Value = 0.006805322652716498
$ws.Range("G6").Value = 13.68414866666667
$ws.Range("H6").Value = 41.052446
$ws.Range("I6").Value = 0.06687192512742757
$ws.Range("J6").Value = 0.06728034099283527
$ws.Range("M6").Value = 80.48302466666667
$ws.Range("N6").Value = 241.449074
$ws.Range("O6").Value = 0.2563744769156197
$ws.Range("P6").Value = 0.2693404563502215
$ws.Range("Q6").Value = 1101.341674681667
$ws.Range("R6").Value = 9912.075072135003
$ws.Range("S6").Value = 0.01714425482488472
$ws.Range("T6").Value = 0.01812131774640877
$ws.Range("I7").Value = 0.622425171752627
$ws.Range("J7").Value = 0.6262265923740385
$ws.Range("M7").Value = 3.029459000000001
$ws.Range("N7").Value = 9.088377000000001
$ws.Range("O7").Value = 0.009650183621689718
$ws.Range("P7").Value = 0.0101382356457613
$ws.Range("Q7").Value = 385.8568236188324
$ws.Range("R7").Value = 3472.711412569492
$ws.Range("S7").Value = 0.006006517198174612
$ws.Range("T7").Value = 0.006348832761130109
$ws.Range("I8").Value = 0.622425171752627
$ws.Range("J8").Value = 0.6262265923740385
$ws.Range("O8").Value = 0.2405828560439687
$ws.Range("P8").Value = 0.2527501840920355
$ws.Range("S8").Value = 0.1497448254939047
$ws.Range("T8").Value = 0.1582788865058663
$ws.Range("I9").Value = 0.622425171752627
$ws.Range("J9").Value = 0.6262265923740385
$ws.Range("M9").Value = 109.55234
$ws.Range("N9").Value = 328.65702
$ws.Range("O9").Value = 0.348973264594696
$ws.Range("P9").Value = 0.3666223700220274
$ws.Range("Q9").Value = 13953.48738253607
$ws.Range("R9").Value = 125581.3864428247
$ws.Range("S9").Value = 0.2172097441524286
$ws.Range("T9").Value = 0.2295886774669881
$ws.Range("I10").Value = 0.622425171752627
$ws.Range("J10").Value = 0.6262265923740385
$ws.Range("M10").Value = 45.33717900000001
$ws.Range("N10").Value = 90.67435800000001
$ws.Range("O10").Value = 0.144419218824026
$ws.Range("P10").Value = 0.1011487538899543
$ws.Range("Q10").Value = 5774.51613663642
$ws.Range("R10").Value = 34647.09681981852
$ws.Range("S10").Value = 0.08989015708092463
$ws.Range("T10").Value = 0.06334203947138636
$ws.Range("I11").Value = 0.622425171752627
$ws.Range("J11").Value = 0.6262265923740385
$ws.Range("M11").Value = 80.48302466666667
$ws.Range("N11").Value = 241.449074
$ws.Range("O11").Value = 0.2563744769156197
$ws.Range("P11").Value = 0.2693404563502215
$ws.Range("Q11").Value = 10250.98020904595
$ws.Range("R11").Value = 92258.82188141355
$ws.Range("S11").Value = 0.1595739278271945
$ws.Range("T11").Value = 0.1686681561686677
$ws.Range("G12").Value = 17.548286
$ws.Range("H12").Value = 52.644858
$ws.Range("I12").Value = 0.08575525566783661
$ws.Range("J12").Value = 0.08627900022715801
$ws.Range("M12").Value = 3.029459000000001
$ws.Range("N12").Value = 9.088377000000001
$ws.Range("O12").Value = 0.009650183621689718
$ws.Range("P12").Value = 0.0101382356457613
$ws.Range("Q12").Value = 53.16181295727402
$ws.Range("R12").Value = 478.4563166154661
$ws.Range("S12").Value = 0.0008275539637195712
$ws.Range("T12").Value = 0.0008747168355836205
$ws.Range("G13").Value = 17.548286
$ws.Range("H13").Value = 52.644858
$ws.Range("I13").Value = 0.08575525566783661
$ws.Range("J13").Value = 0.08627900022715801
$ws.Range("O13").Value = 0.2405828560439687
$ws.Range("P13").Value = 0.2527501840920355
$ws.Range("Q13").Value = 1325.344811573314
$ws.Range("R13").Value = 11928.10330415983
$ws.Range("S13").Value = 0.02063124432934886
$ws.Range("T13").Value = 0.02180703319069096
$ws.Range("G14").Value = 17.548286
$ws.Range("H14").Value = 52.644858
$ws.Range("I14").Value = 0.08575525566783661
$ws.Range("J14").Value = 0.08627900022715801
$ws.Range("M14").Value = 109.55234
$ws.Range("N14").Value = 328.65702
$ws.Range("O14").Value = 0.348973264594696
$ws.Range("P14").Value = 0.3666223700220274
$ws.Range("Q14").Value = 1922.45579428924
$ws.Range("R14").Value = 17302.10214860316
$ws.Range("S14").Value = 0.02992629152655775
$ws.Range("T14").Value = 0.03163181154641171
$ws.Range("G15").Value = 17.548286
$ws.Range("H15").Value = 52.644858
$ws.Range("I15").Value = 0.08575525566783661
$ws.Range("J15").Value = 0.08627900022715801
$ws.Range("M15").Value = 45.33717900000001
$ws.Range("N15").Value = 90.67435800000001
$ws.Range("O15").Value = 0.144419218824026
$ws.Range("P15").Value = 0.1011487538899543
$ws.Range("Q15").Value = 795.5897835251941
$ws.Range("R15").Value = 4773.538701151165
$ws.Range("S15").Value = 0.01238470703360359
$ws.Range("T15").Value = 0.008727013359848117
$ws.Range("G16").Value = 17.548286
$ws.Range("H16").Value = 52.644858
$ws.Range("I16").Value = 0.08575525566783661
$ws.Range("J16").Value = 0.08627900022715801
$ws.Range("M16").Value = 80.48302466666667
$ws.Range("N16").Value = 241.449074
$ws.Range("O16").Value = 0.2563744769156197
$ws.Range("P16").Value = 0.2693404563502215
$ws.Range("Q16").Value = 1412.339134995721
$ws.Range("R16").Value = 12711.05221496149
$ws.Range("S16").Value = 0.02198545881460684
$ws.Range("T16").Value = 0.0232384252946236
$ws.Range("G17").Value = 3.726573
$ws.Range("H17").Value = 7.453145999999999
$ws.Range("I17").Value = 0.01821107887003078
$ws.Range("J17").Value = 0.01221486788751604
$ws.Range("M17").Value = 3.029459000000001
$ws.Range("N17").Value = 9.088377000000001
$ws.Range("O17").Value = 0.009650183621689718
$ws.Range("P17").Value = 0.0101382356457613
$ws.Range("Q17").Value = 11.289500114007
$ws.Range("R17").Value = 67.737000684042
$ws.Range("S17").Value = 0.0001757402550448707
$ws.Range("T17").Value = 0.0001238372090254801
$ws.Range("G18").Value = 3.726573
$ws.Range("H18").Value = 7.453145999999999
$ws.Range("I18").Value = 0.01821107887003078
$ws.Range("J18").Value = 0.01221486788751604
$ws.Range("O18").Value = 0.2405828560439687
$ws.Range("P18").Value = 0.2527501840920355
$ws.Range("Q18").Value = 281.451658042227
$ws.Range("R18").Value = 1688.709948253362
$ws.Range("S18").Value = 0.004381273366193973
$ws.Range("T18").Value = 0.003087310107229571
$ws.Range("G19").Value = 3.726573
$ws.Range("H19").Value = 7.453145999999999
$ws.Range("I19").Value = 0.01821107887003078
$ws.Range("J19").Value = 0.01221486788751604
$ws.Range("M19").Value = 109.55234
$ws.Range("N19").Value = 328.65702
$ws.Range("O19").Value = 0.348973264594696
$ws.Range("P19").Value = 0.3666223700220274
$ws.Range("Q19").Value = 408.25479233082
$ws.Range("R19").Value = 2449.52875398492
$ws.Range("S19").Value = 0.006355179645066127
$ws.Range("T19").Value = 0.004478243814427085
$ws.Range("G20").Value = 3.726573
$ws.Range("H20").Value = 7.453145999999999
$ws.Range("I20").Value = 0.01821107887003078
$ws.Range("J20").Value = 0.01221486788751604
$ws.Range("M20").Value = 45.33717900000001
$ws.Range("N20").Value = 90.67435800000001
$ws.Range("O20").Value = 0.144419218824026
$ws.Range("P20").Value = 0.1011487538899543
$ws.Range("Q20").Value = 168.952307157567
$ws.Range("R20").Value = 675.8092286302681
$ws.Range("S20").Value = 0.002630029784352571
$ws.Range("T20").Value = 0.001235518665752666
$ws.Range("G21").Value = 3.726573
$ws.Range("H21").Value = 7.453145999999999
$ws.Range("I21").Value = 0.01821107887003078
$ws.Range("J21").Value = 0.01221486788751604
$ws.Range("M21").Value = 80.48302466666667
$ws.Range("N21").Value = 241.449074
$ws.Range("O21").Value = 0.2563744769156197
$ws.Range("P21").Value = 0.2693404563502215
$ws.Range("Q21").Value = 299.925866681134
$ws.Range("R21").Value = 1799.555200086804
$ws.Range("S21").Value = 0.004668855819373234
$ws.Range("T21").Value = 0.003289958091081236
$ws.Range("G22").Value = 42.30495733333333
$ws.Range("H22").Value = 126.914872
$ws.Range("I22").Value = 0.206736568582078
$ws.Range("J22").Value = 0.2079991985184523
$ws.Range("M22").Value = 3.029459000000001
$ws.Range("N22").Value = 9.088377000000001
$ws.Range("O22").Value = 0.009650183621689718
$ws.Range("P22").Value = 0.0101382356457613
$ws.Range("Q22").Value = 128.1611337380827
$ws.Range("R22").Value = 1153.450203642744
$ws.Range("S22").Value = 0.001995045848135102
$ws.Range("T22").Value = 0.002108744888709554
$ws.Range("G23").Value = 42.30495733333333
$ws.Range("H23").Value = 126.914872
$ws.Range("I23").Value = 0.206736568582078
$ws.Range("J23").Value = 0.2079991985184523
$ws.Range("O23").Value = 0.2405828560439687
$ws.Range("P23").Value = 0.2527501840920355
$ws.Range("Q23").Value = 3195.107243269443
$ws.Range("R23").Value = 28755.96518942498
$ws.Range("S23").Value = 0.04973727411820613
$ws.Range("T23").Value = 0.05257183571653466
$ws.Range("G24").Value = 42.30495733333333
$ws.Range("H24").Value = 126.914872
$ws.Range("I24").Value = 0.206736568582078
$ws.Range("J24").Value = 0.2079991985184523
$ws.Range("M24").Value = 109.55234
$ws.Range("N24").Value = 328.65702
$ws.Range("O24").Value = 0.348973264594696
$ws.Range("P24").Value = 0.3666223700220274
$ws.Range("Q24").Value = 4634.607069466827
$ws.Range("R24").Value = 41711.46362520144
$ws.Range("S24").Value = 0.07214553524919302
$ws.Range("T24").Value = 0.07625715912351715
$ws.Range("G25").Value = 42.30495733333333
$ws.Range("H25").Value = 126.914872
$ws.Range("I25").Value = 0.206736568582078
$ws.Range("J25").Value = 0.2079991985184523
$ws.Range("M25").Value = 45.33717900000001
$ws.Range("N25").Value = 90.67435800000001
$ws.Range("O25").Value = 0.144419218824026
$ws.Range("P25").Value = 0.1011487538899543
$ws.Range("Q25").Value = 1917.987423208696
$ws.Range("R25").Value = 11507.92453925218
$ws.Range("S25").Value = 0.02985673373698339
$ws.Range("T25").Value = 0.02103885974025068
$ws.Range("G26").Value = 42.30495733333333
$ws.Range("H26").Value = 126.914872
$ws.Range("I26").Value = 0.206736568582078
$ws.Range("J26").Value = 0.2079991985184523
$ws.Range("M26").Value = 80.48302466666667
$ws.Range("N26").Value = 241.449074
$ws.Range("O26").Value = 0.2563744769156197
$ws.Range("P26").Value = 0.2693404563502215
$ws.Range("Q26").Value = 3404.830924580948
$ws.Range("R26").Value = 30643.47832122853
$ws.Range("S26").Value = 0.05300197962956037
$ws.Range("T26").Value = 0.05602259904944026
